# Generate Report for Handoff
#
# The localization status report moves from "In Translation" to
# "Ready for handoff": the per-language status cells on the Overview
# sheet and the Status column on each language sheet are updated, the
# "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# are refreshed, and the Status column is widened on every sheet to fit
# the new, longer label.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps (kept as text, matching the source) -
$overview.Range("G2").Value = "2016-08-22 04:49:33"
$zhcn.Range("H2").Value     = "2016-08-22 04:49:29"
$dede.Range("H2").Value     = "2016-08-22 04:49:33"

# Re-assert the existing "yyyy-mm-dd HH:mm:ss" display format on the
# datetime cells so it stays attached after the value update (K2 is
# untouched by this change but carries the same format and is refreshed
# defensively so a save round-trip doesn't drop it).
$overview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("H2").NumberFormat     = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("K2").NumberFormat     = "yyyy-mm-dd HH:mm:ss"
$dede.Range("H2").NumberFormat     = "yyyy-mm-dd HH:mm:ss"
$dede.Range("K2").NumberFormat     = "yyyy-mm-dd HH:mm:ss"

# --- Widen the Status column so the longer text fits -----------------
# (16.3333.. is the character width whose pixel-grid snap lands closest
# to the authored 17.216 stored column width.)
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth     = 16.3333333333333
